$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add "Passed" results for Save Game / Load Game rows (20-21) ---
$ws.Range("C20").Value = "Passed"
$ws.Range("E20").Value = "When a user clicks save, data is saved to a file for a user to load later"

$ws.Range("C21").Value = "Passed"
$ws.Range("E21").Value = "when you click load, a saved game is loaded into the gameboard to continue"

# --- Add new rows 23-24 documenting the remaining Failed side-pin issue ---
$ws.Range("A23").Value = "Save Game"
$ws.Range("C23").Value = "Failed"
$ws.Range("E23").Value = "When a user clicks save the data is put into arrays however a file isnt created"

$ws.Range("A24").Value = "Load Game"
$ws.Range("C24").Value = "Failed"
$ws.Range("E24").Value = "as a file isnt created a past game cannot load to gameboard"

# --- Widen column E slightly to fit the new text ---
$ws.Columns.Item(5).ColumnWidth = 64.3

# --- Move the active selection to reflect where the author left off editing ---
$ws.Range("F23").Select()
